$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (kept as Text via leading apostrophe, like the
# original cells which are stored as text/inline strings, not numbers).
$updates = @{
    "D2" = "308.94"
    "E2" = "-2.27%"
    "D3" = "48.65"
    "E3" = "7.37%"
    "E4" = "0.81%"
    "D5" = "0.07740"
    "E5" = "-4.28%"
    "D6" = "4.510"
    "E6" = "-0.58%"
    "D7" = "1.297"
    "E7" = "18.71%"
    "D8" = "1.555"
    "E8" = "-7.20%"
    "D9" = "0.1232"
    "E9" = "-5.63%"
    "D10" = "0.1918"
    "E10" = "-0.80%"
    "D11" = "0.09194"
    "E11" = "-2.71%"
    "D12" = "0.04555"
    "E12" = "7.37%"
    "D13" = "0.1048"
    "E13" = "0.38%"
    "E14" = "-3.68%"
    "D15" = "0.04199"
    "E15" = "-1.50%"
    "D16" = "0.005882"
    "E16" = "0.04%"
    "D17" = "3.346"
    "E17" = "-1.57%"
    "D18" = "2.400"
    "E18" = "-0.43%"
    "D19" = "0.3439"
    "E19" = "2.00%"
    "D20" = "8.115"
    "E20" = "-1.37%"
    "D21" = "0.1365"
    "E21" = "-1.49%"
    "D22" = "0.3037"
    "E22" = "-3.45%"
    "D23" = "0.001299"
    "E23" = "1.31%"
    "D24" = "0.004093"
    "E24" = "-2.94%"
    "D25" = "0.0001362"
    "E25" = "1.26%"
    "D38" = "0.02569"
    "E38" = "-5.04%"
    "D39" = "0.05736"
    "E39" = "4.94%"
    "D40" = "0.009237"
    "E40" = "57.52%"
    "D41" = "0.007979"
    "E41" = "2.88%"
    "D42" = "0.1420"
    "E42" = "-0.13%"
    "D43" = "0.008391"
    "E43" = "13.81%"
    "D44" = "0.007772"
    "E44" = "-9.48%"
    "D45" = "0.3385"
    "E45" = "7.83%"
    "D46" = "0.00006916"
    "E46" = "1.69%"
    "E47" = "1.17%"
    "D48" = "0.05509"
    "E48" = "-20.69%"
    "E49" = "1.19%"
    "E50" = "1.17%"
    "E51" = "1.17%"
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.Value = "'" + $updates[$cell]
    $range.Style = "Normal"
}

